# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must stay stored as TEXT
# (matches the source data which keeps these as plain strings).
$textCells = @("D4", "D5", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '30.804.48'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = '1.938.68'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '243.83'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("D7").Value = '0.4887'
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").Value = '  -0.44%  '
$ws.Range("D9").Value = '0.06894'
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("D10").Value = '19.31'
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("D11").Value = '105.09'
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("D12").Value = '0.07791'
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("D13").Value = '1.937.12'
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("D14").Value = '5.355'
$ws.Range("E14").Value = '  -2.18%  '
$ws.Range("D15").Value = '0.7026'
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").Value = '273.13'
$ws.Range("E16").Value = '  -3.04%  '
$ws.Range("D17").Value = '30.807.25'
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").Value = '0.000007735'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '5.664'
$ws.Range("E19").Value = '  +3.12%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '13.16'
$ws.Range("E20").Value = '  -0.59%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '6.538'
$ws.Range("E23").Value = '  +0.63%  '
$ws.Range("D24").Value = '9.819'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = '165.17'
$ws.Range("E25").Value = '  -2.29%  '
$ws.Range("D26").Value = '19.61'
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("D27").Value = '2.164'
$ws.Range("E27").Value = '  -2.26%  '
$ws.Range("D28").Value = '0.1038'
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("D29").Value = '1.386'
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("D30").Value = '4.644'
$ws.Range("E30").Value = '  +1.75%  '
$ws.Range("D31").Value = '1.561'
$ws.Range("E31").Value = '  -1.40%  '
$ws.Range("D32").Value = '4.432'
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("D34").Value = '0.7604'
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("D35").Value = '1.153'
$ws.Range("E35").Value = '  -1.69%  '
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("E37").Value = '  +0.58%  '
$ws.Range("D38").Value = '0.02012'
$ws.Range("E38").Value = '  -0.59%  '
$ws.Range("D39").Value = '79.58'
$ws.Range("E39").Value = '  +6.17%  '
$ws.Range("D40").Value = '2.669'
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("D41").Value = '6.509'
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").Value = '2.087'
$ws.Range("E42").Value = '  -3.24%  '
$ws.Range("D43").Value = '0.9065'
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("D44").Value = '0.4461'
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = '108.64'
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").Value = '7.879'
$ws.Range("E46").Value = '  -3.53%  '
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").Value = '997.10'
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("D49").Value = '0.1252'
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("D51").Value = '9.232'
$ws.Range("E51").Value = '  -1.88%  '

# Restore default (Normal) style on the forced-text cells so we don't
# leave a stray style index behind (matches original formatting).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

